$wb = $excel.ActiveWorkbook

# --- Hoja "Clientes": se inserta un nuevo registro y se reordenan filas ---
$wsClientes = $wb.Worksheets.Item("Clientes")

# Nueva fila 2: registro de "Ester"
$wsClientes.Cells.Item(2, 1).Value = 43870330
$wsClientes.Cells.Item(2, 2).Value = "Ester"
$wsClientes.Cells.Item(2, 3).Value = 3105206465

# Fila 3 ("Carlos"): la cedula pasa a ser numerica y se corrige el telefono
$wsClientes.Cells.Item(3, 1).Value = 71722939
$wsClientes.Cells.Item(3, 2).Value = "Carlos"
$wsClientes.Cells.Item(3, 3).Value = 3022350912

# Fila 4 ("Alejandro"): registro existente se reubica al final
$wsClientes.Cells.Item(4, 1).Value = 1013337950
$wsClientes.Cells.Item(4, 2).Value = "Alejandro"
$wsClientes.Cells.Item(4, 3).Value = 3015305600

# --- Hoja "Productos" ---
$wsProductos = $wb.Worksheets.Item("Productos")

# Ajustar ancho de la columna C (Marca) -> el valor stored (OOXML) resultante es 8
# (Excel aplica un relleno ~5/6 al convertir "caracteres" a unidades almacenadas)
$wsProductos.Columns.Item(3).ColumnWidth = 7.15

# Fila 6: "Papas" - se actualiza la marca y los precios
$wsProductos.Cells.Item(6, 3).Value = "Exito"
$wsProductos.Cells.Item(6, 4).Value = 30000
$wsProductos.Cells.Item(6, 5).Value = 20000

# Fila 7: "Frijoles" - se actualiza la marca y los precios
$wsProductos.Cells.Item(7, 3).Value = "Zenu"
$wsProductos.Cells.Item(7, 4).Value = 1500
$wsProductos.Cells.Item(7, 5).Value = 3500

# Nuevos productos agregados
$wsProductos.Cells.Item(8, 1).Value = "Zapatos"
$wsProductos.Cells.Item(8, 2).Value = 7702003141516
$wsProductos.Cells.Item(8, 3).Value = "Adidas"
$wsProductos.Cells.Item(8, 4).Value = 20
$wsProductos.Cells.Item(8, 5).Value = 30
$wsProductos.Cells.Item(8, 6).Value = 20
$wsProductos.Cells.Item(8, 7).Value = $true
$wsProductos.Cells.Item(8, 8).Value = "04/06/2024 16:19"

$wsProductos.Cells.Item(9, 1).Value = "Camiseta"
$wsProductos.Cells.Item(9, 2).Value = 7702034121618
$wsProductos.Cells.Item(9, 3).Value = "Rifle"
$wsProductos.Cells.Item(9, 4).Value = 10
$wsProductos.Cells.Item(9, 5).Value = 50
$wsProductos.Cells.Item(9, 6).Value = 15
$wsProductos.Cells.Item(9, 7).Value = $true
$wsProductos.Cells.Item(9, 8).Value = "04/06/2024 17:03"

$wsProductos.Cells.Item(10, 1).Value = "Gorra"
$wsProductos.Cells.Item(10, 2).Value = 7703001565152
$wsProductos.Cells.Item(10, 3).Value = "Calvin"
$wsProductos.Cells.Item(10, 4).Value = 40
$wsProductos.Cells.Item(10, 5).Value = 70
$wsProductos.Cells.Item(10, 6).Value = 20
$wsProductos.Cells.Item(10, 7).Value = $true
$wsProductos.Cells.Item(10, 8).Value = "04/06/2024 17:26"
